$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "CMS ID"
$ws.Range("B1").Value = "Title"
$ws.Range("C1").Value = "PUBLISHING PLACE"
$ws.Range("D1").Value = "Publisher"
$ws.Range("E1").Value = "Language"
$ws.Range("F1").Value = "Comments"

# --- Row 2 ---
$ws.Range("A2").Value = 1234
$ws.Range("B2").Value = "Göttinger Tageblatt"
$ws.Range("C2").Value = "Göttingen"
$ws.Range("D2").Value = "Verlag des Wissens"
$ws.Range("E2").Value = "de"
$ws.Range("F2").Value = "Göttinger Stadtzeitung"

# --- Row 3 ---
$ws.Range("A3").Value = 5678
$ws.Range("B3").Value = "Rhein-Hunsrück-Zeitung"
$ws.Range("C3").Value = "Simmern"
$ws.Range("D3").Value = "Rhein-Zeitung"
$ws.Range("E3").Value = "de"
$ws.Range("F3").Value = "Zeitung des Landkreises Rhein-Hunsrück"

# --- Row 4 ---
$ws.Range("A4").Value = 9101
$ws.Range("B4").Value = "Rosdorfer Gemeindeblatt"
$ws.Range("C4").Value = "Rosdorf"
$ws.Range("D4").Value = "Gemeinde Rosdorf"
$ws.Range("E4").Value = "de"
$ws.Range("F4").Value = "Kleines Gemeindeblättchen"

# --- Selection, matching final cursor position recorded in the file ---
$ws.Range("F5").Select() | Out-Null
